$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, $Text) {
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '64.312.99'
Set-TextValue $ws.Range("E2") '  +0.95%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.490.94'
Set-TextValue $ws.Range("E3") '  -0.02%  '

# Row 4
Set-TextValue $ws.Range("E4") '  -0.02%  '

# Row 5
Set-TextValue $ws.Range("D5") '586.44'
Set-TextValue $ws.Range("E5") '  +0.35%  '

# Row 6
Set-TextValue $ws.Range("D6") '134.63'
Set-TextValue $ws.Range("E6") '  +2.57%  '

# Row 8
Set-TextValue $ws.Range("E8") '  -0.72%  '

# Row 9
Set-TextValue $ws.Range("E9") '  +0.85%  '

# Row 10
Set-TextValue $ws.Range("E10") '  +0.77%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.378'
Set-TextValue $ws.Range("E11") '  -1.28%  '

# Row 12
Set-TextValue $ws.Range("D12") '4.085.29'
Set-TextValue $ws.Range("E12") '  +0.11%  '

# Row 13
Set-TextValue $ws.Range("E13") '  +1.30%  '

# Row 14
Set-TextValue $ws.Range("E14") '  +1.21%  '

# Row 15
Set-TextValue $ws.Range("D15") '3.489.38'
Set-TextValue $ws.Range("E15") '  -0.17%  '

# Row 16
Set-TextValue $ws.Range("D16") '64.311.49'
Set-TextValue $ws.Range("E16") '  +0.71%  '

# Row 17
Set-TextValue $ws.Range("D17") '25.66'
Set-TextValue $ws.Range("E17") '  -6.97%  '

# Row 18
Set-TextValue $ws.Range("D18") '9.89'
Set-TextValue $ws.Range("E18") '  -2.71%  '

# Row 19
Set-TextValue $ws.Range("D19") '5.75'
Set-TextValue $ws.Range("E19") '  +1.57%  '

# Row 20
Set-TextValue $ws.Range("D20") '13.63'
Set-TextValue $ws.Range("E20") '  -4.94%  '

# Row 21
Set-TextValue $ws.Range("D21") '386.84'
Set-TextValue $ws.Range("E21") '  +0.18%  '

# Row 22
Set-TextValue $ws.Range("E22") '  -1.80%  '

# Row 23
Set-TextValue $ws.Range("D23") '3.629.80'
Set-TextValue $ws.Range("E23") '  -0.02%  '

# Row 24
Set-TextValue $ws.Range("D24") '74.33'
Set-TextValue $ws.Range("E24") '  +1.83%  '

# Row 26
Set-TextValue $ws.Range("D26") '5.71'
Set-TextValue $ws.Range("E26") '  -0.42%  '

# Row 27
Set-TextValue $ws.Range("D27") '0.0000113'
Set-TextValue $ws.Range("E27") '  +0.59%  '

# Row 28
Set-TextValue $ws.Range("E28") '  +0.11%  '

# Row 29
Set-TextValue $ws.Range("D29") '7.40'
Set-TextValue $ws.Range("E29") '  -0.13%  '

# Row 30
Set-TextValue $ws.Range("D30") '1.49'
Set-TextValue $ws.Range("E30") '  -5.24%  '

# Row 31
Set-TextValue $ws.Range("B31") 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range("C31") 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D31") '8.26'
Set-TextValue $ws.Range("E31") '  +0.68%  '

# Row 32
Set-TextValue $ws.Range("B32") 'PancakeSwap'
Set-TextValue $ws.Range("C32") 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D32") '2.23'
Set-TextValue $ws.Range("E32") '  -0.63%  '

# Row 33
Set-TextValue $ws.Range("D33") '3.511.50'
Set-TextValue $ws.Range("E33") '  +0.56%  '

# Row 34
Set-TextValue $ws.Range("E34") '  -0.01%  '

# Row 35
Set-TextValue $ws.Range("E35") '  +3.75%  '

# Row 36
Set-TextValue $ws.Range("E36") '  -1.09%  '

# Row 37
Set-TextValue $ws.Range("E37") '  -0.04%  '

# Row 38
Set-TextValue $ws.Range("E38") '  +0.09%  '

# Row 39
Set-TextValue $ws.Range("D39") '1.54'
Set-TextValue $ws.Range("E39") '  -1.39%  '

# Row 40
Set-TextValue $ws.Range("D40") '162.87'
Set-TextValue $ws.Range("E40") '  -2.64%  '

# Row 41
Set-TextValue $ws.Range("D41") '0.0781'
Set-TextValue $ws.Range("E41") '  -2.31%  '

# Row 42
Set-TextValue $ws.Range("D42") '0.804'
Set-TextValue $ws.Range("E42") '  -0.80%  '

# Row 43
Set-TextValue $ws.Range("B43") 'EnergySwap'
Set-TextValue $ws.Range("C43") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D43") '25.50'
Set-TextValue $ws.Range("E43") '  -5.57%  '

# Row 44
Set-TextValue $ws.Range("B44") 'FirstDigitalUSD'
Set-TextValue $ws.Range("C44") 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range("D44") '1.00'
Set-TextValue $ws.Range("E44") '  -0.03%  '

# Row 45
Set-TextValue $ws.Range("E45") '  +0.58%  '

# Row 46
Set-TextValue $ws.Range("E46") '  +0.50%  '

# Row 47
Set-TextValue $ws.Range("E47") '  +1.77%  '

# Row 48
Set-TextValue $ws.Range("E48") '  -2.41%  '

# Row 49
Set-TextValue $ws.Range("D49") '2.474.46'
Set-TextValue $ws.Range("E49") '  +1.77%  '

# Row 50
Set-TextValue $ws.Range("D50") '6.74'
Set-TextValue $ws.Range("E50") '  -1.70%  '

# Row 51
Set-TextValue $ws.Range("E51") '  +1.39%  '
